$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value "Mkdad" to cell A7
$ws.Range("A7").Value = "Mkdad"

# Update the selection to A7
$ws.Range("A7").Select()
